# Azure Network Documentation Template.xlsx
# Add the missing "AzureBastion" subnet row (row 3, columns D:G) that
# documents the Bastion subnet: its address space, the service deployed
# there, and whether an NSG / UDR is attached.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# Fill in the new subnet-documentation cells for row 3 (AzureBastion).
$ws.Range("D3").Value = "10.0.1.0/26"
$ws.Range("E3").Value = "Bastion"
$ws.Range("F3").Value = "Y"
$ws.Range("G3").Value = "N"

# Match the centered alignment used by the NSG/UDR columns in row 2.
$ws.Range("F3:G3").HorizontalAlignment = $xlCenter

# Leave the selection on the newly-documented Bastion service cell.
$ws.Range("E3").Select()
